$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="30.010.60"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value2 = '  -0.42%  '

$ws.Range("D3").Formula = '="1.873.57"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value2 = '  -2.57%  '

$ws.Range("D4").Formula = '="1.0000"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value2 = '  +0.01%  '

$ws.Range("D5").Formula = '="319.67"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)

$ws.Range("D6").Formula = '="0.9997"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)

$ws.Range("D7").Formula = '="0.5097"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value2 = '  -3.06%  '

$ws.Range("D8").Formula = '="0.3957"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value2 = '  -2.64%  '

$ws.Range("D9").Formula = '="0.08217"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value2 = '  -3.93%  '

$ws.Range("D10").Formula = '="42.17"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value2 = '  -2.61%  '

$ws.Range("D11").Formula = '="1.096"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value2 = '  -3.11%  '

$ws.Range("D12").Formula = '="23.91"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value2 = '  +5.59%  '

$ws.Range("D13").Formula = '="1.866.71"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value2 = '  -2.92%  '

$ws.Range("D14").Formula = '="6.312"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value2 = '  -1.93%  '

$ws.Range("D15").Formula = '="7.211"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value2 = '  -2.78%  '

$ws.Range("D16").Formula = '="0.9998"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value2 = '  +0.00%  '

$ws.Range("D17").Formula = '="92.09"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value2 = '  -4.80%  '

$ws.Range("E18").Value2 = '  -2.84%  '

$ws.Range("D19").Formula = '="0.06391"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value2 = '  -4.89%  '

$ws.Range("D20").Formula = '="18.06"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value2 = '  -1.43%  '

$ws.Range("D21").Formula = '="0.9999"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value2 = '  +0.05%  '

$ws.Range("D22").Formula = '="29.993.42"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value2 = '  -0.45%  '

$ws.Range("D23").Formula = '="5.848"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value2 = '  -3.57%  '

$ws.Range("D24").Formula = '="11.15"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value2 = '  -1.44%  '

$ws.Range("E25").Value2 = '  -2.24%  '

$ws.Range("D26").Formula = '="2.085.08"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value2 = '  -2.69%  '

$ws.Range("D27").Formula = '="160.73"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value2 = '  +0.25%  '

$ws.Range("D28").Formula = '="21.14"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value2 = '  -0.38%  '

$ws.Range("D29").Formula = '="2.247"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value2 = '  -9.38%  '

$ws.Range("E30").Value2 = '  -1.51%  '

$ws.Range("D31").Formula = '="1.073"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value2 = '  -0.98%  '

$ws.Range("D32").Formula = '="0.1036"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value2 = '  -2.26%  '

$ws.Range("D33").Formula = '="5.949"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value2 = '  -3.07%  '

$ws.Range("D34").Formula = '="3.714"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value2 = '  +1.86%  '

$ws.Range("D35").Formula = '="0.02442"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value2 = '  -3.39%  '

$ws.Range("D36").Formula = '="5.236"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value2 = '  +0.26%  '

$ws.Range("D37").Formula = '="0.06381"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value2 = '  -3.47%  '

$ws.Range("D38").Formula = '="0.2146"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value2 = '  -3.82%  '

$ws.Range("D39").Formula = '="1.180"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value2 = '  -4.81%  '

$ws.Range("D40").Formula = '="8.554"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value2 = '  -6.07%  '

$ws.Range("B41").Value2 = 'Aptos'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Formula = '="11.42"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value2 = '  -2.60%  '

$ws.Range("B42").Value2 = 'TheSandbox'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Formula = '="0.6324"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value2 = '  -3.95%  '

$ws.Range("D43").Formula = '="1.205"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value2 = '  -3.15%  '

$ws.Range("D44").Formula = '="0.9985"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value2 = '  +0.01%  '

$ws.Range("D45").Formula = '="13.01"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value2 = '  -2.21%  '

$ws.Range("D46").Formula = '="0.5919"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value2 = '  -4.68%  '

$ws.Range("E47").Value2 = '  -3.92%  '

$ws.Range("D48").Formula = '="2.020"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value2 = '  -3.66%  '

$ws.Range("D49").Formula = '="122.76"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value2 = '  -2.04%  '

$ws.Range("D50").Formula = '="1.207"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value2 = '  -3.43%  '

$ws.Range("D51").Formula = '="1.123"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value2 = '  -2.84%  '

$excel.CutCopyMode = $false
